# Apply cryptos list update (cell-level changes per diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.373.67"
$ws.Cells.Item(2, 5).Value = "  +0.43%  "
$ws.Cells.Item(3, 4).Value = "1.612.33"
$ws.Cells.Item(3, 5).Value = "  +1.48%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.998"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "213.73"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.502"
$ws.Cells.Item(6, 5).Value = "  +0.31%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(9, 5).Value = "  +0.22%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.16"
$ws.Cells.Item(10, 5).Value = "  -0.96%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0854"
$ws.Cells.Item(11, 5).Value = "  +0.44%  "
$ws.Cells.Item(12, 4).Value = "1.838.19"
$ws.Cells.Item(12, 5).Value = "  +1.42%  "
$ws.Cells.Item(13, 4).Value = "1.630.00"
$ws.Cells.Item(13, 5).Value = "  +2.43%  "
$ws.Cells.Item(14, 5).Value = "  +0.08%  "
$ws.Cells.Item(15, 5).Value = "  -1.55%  "
$ws.Cells.Item(16, 5).Value = "  +0.41%  "
$ws.Cells.Item(17, 4).Value = "26.379.55"
$ws.Cells.Item(17, 5).Value = "  +0.48%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0725"
$ws.Cells.Item(18, 5).Value = "  -0.12%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "220.80"
$ws.Cells.Item(19, 5).Value = "  +3.81%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.54"
$ws.Cells.Item(20, 5).Value = "  +1.88%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.00"
$ws.Cells.Item(21, 5).Value = "  -0.13%  "
$ws.Cells.Item(22, 5).Value = "  +1.73%  "
$ws.Cells.Item(23, 5).Value = "  +0.72%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.14"
$ws.Cells.Item(24, 5).Value = "  -0.10%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.93"
$ws.Cells.Item(25, 5).Value = "  +0.66%  "
$ws.Cells.Item(26, 5).Value = "  -0.14%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.03"
$ws.Cells.Item(27, 5).Value = "  -0.48%  "
$ws.Cells.Item(28, 5).Value = "  +1.37%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.27"
$ws.Cells.Item(29, 5).Value = "  +0.52%  "
$ws.Cells.Item(30, 5).Value = "  -0.11%  "
$ws.Cells.Item(31, 5).Value = "  +0.31%  "
$ws.Cells.Item(32, 5).Value = "  +0.43%  "
$ws.Cells.Item(33, 4).Value = "1.446.72"
$ws.Cells.Item(33, 5).Value = "  +8.39%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.98"
$ws.Cells.Item(34, 5).Value = "  +1.46%  "
$ws.Cells.Item(35, 5).Value = "  -0.82%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.48"
$ws.Cells.Item(36, 5).Value = "  +0.22%  "
$ws.Cells.Item(37, 5).Value = "  -5.01%  "
$ws.Cells.Item(38, 5).Value = "  -0.09%  "
$ws.Cells.Item(40, 5).Value = "  +1.90%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.999"
$ws.Cells.Item(41, 5).Value = "  -0.15%  "
$ws.Cells.Item(42, 5).Value = "  +1.90%  "
$ws.Cells.Item(43, 4).Value = "1.750.75"
$ws.Cells.Item(43, 5).Value = "  +1.51%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.762"
$ws.Cells.Item(44, 5).Value = "  -0.01%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "61.66"
$ws.Cells.Item(45, 5).Value = "  -0.30%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.906"
$ws.Cells.Item(46, 5).Value = "  -11.72%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "88.07"
$ws.Cells.Item(47, 5).Value = "  +2.75%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0105"
$ws.Cells.Item(48, 5).Value = "  +0.93%  "
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.49"
$ws.Cells.Item(49, 5).Value = "  +0.33%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0502"
$ws.Cells.Item(50, 5).Value = "  -0.04%  "
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0961"
$ws.Cells.Item(51, 5).Value = "  -1.29%  "
